$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows of admin-login test data (username/aa, password/ap)
$ws.Range("A4").Value = "username"
$ws.Range("B4").Value = "aa"
$ws.Range("A5").Value = "password"
$ws.Range("B5").Value = "ap"

# Match the highlighted (yellow fill) text-formatted style used for the new rows
$ws.Range("A4:B5").Interior.Color = 65535

# Update the selected cell as recorded in the saved workbook
$ws.Range("L7").Select()
